$wb = $excel.ActiveWorkbook

# --- 1. Update the "time_taken" timestamps (column F) on the existing "data" sheet ---
$ws1 = $wb.Worksheets.Item("data")

$ws1.Range("F2").Value  = "2021-10-05 14:35:18.618328"
$ws1.Range("F3").Value  = "2021-10-05 14:35:18.618336"
$ws1.Range("F4").Value  = "2021-10-05 14:35:18.618340"
$ws1.Range("F5").Value  = "2021-10-05 14:35:18.618342"
$ws1.Range("F6").Value  = "2021-10-05 14:35:18.618345"
$ws1.Range("F7").Value  = "2021-10-05 14:35:18.618348"
$ws1.Range("F8").Value  = "2021-10-05 14:35:18.618351"
$ws1.Range("F9").Value  = "2021-10-05 14:35:18.618353"
$ws1.Range("F10").Value = "2021-10-05 14:35:18.618356"
$ws1.Range("F11").Value = "2021-10-05 14:35:18.618359"
$ws1.Range("F12").Value = "2021-10-05 14:35:18.618361"
$ws1.Range("F13").Value = "2021-10-05 14:35:18.618364"
$ws1.Range("F14").Value = "2021-10-05 14:35:18.618367"
$ws1.Range("F15").Value = "2021-10-05 14:35:18.618369"
$ws1.Range("F16").Value = "2021-10-05 14:35:18.618372"
$ws1.Range("F17").Value = "2021-10-05 14:35:18.618375"
$ws1.Range("F18").Value = "2021-10-05 14:35:18.618377"
$ws1.Range("F19").Value = "2021-10-05 14:35:18.618380"
$ws1.Range("F20").Value = "2021-10-05 14:35:18.618383"
$ws1.Range("F21").Value = "2021-10-05 14:35:18.618385"
$ws1.Range("F22").Value = "2021-10-05 14:35:18.618388"
$ws1.Range("F23").Value = "2021-10-05 14:35:18.618391"
$ws1.Range("F24").Value = "2021-10-05 14:35:18.618394"
$ws1.Range("F25").Value = "2021-10-05 14:35:18.618396"
$ws1.Range("F26").Value = "2021-10-05 14:35:18.618399"
$ws1.Range("F27").Value = "2021-10-05 14:35:18.618402"
$ws1.Range("F28").Value = "2021-10-05 14:35:18.618404"
$ws1.Range("F29").Value = "2021-10-05 14:35:18.618407"
$ws1.Range("F30").Value = "2021-10-05 14:35:18.618410"
$ws1.Range("F31").Value = "2021-10-05 14:35:18.618412"
$ws1.Range("F32").Value = "2021-10-05 14:35:18.618415"
$ws1.Range("F33").Value = "2021-10-05 14:35:18.618418"

# --- 2. Add a new "metadata" sheet after "data", describing the panel pull itself ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "metadata"

# Re-use the bold/bordered header style and the row-index style from "data"
# (copy formats only, so the same style slot gets reused instead of a new one
# being allocated for every cell).
$ws1.Range("B1").Copy() | Out-Null
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Progressive Myoclonic Epilepsy"
$ws2.Range("C2").Value = 331
# "data_version" is a literal text value ("0.11"), not a number. A plain
# Value assignment would be auto-coerced to a numeric cell, so build it as a
# text formula and flatten the formula back down to a literal value/string.
$ws2.Range("D2").Formula = '="0.11"'
$ws2.Range("D2").Copy() | Out-Null
$ws2.Range("D2").PasteSpecial(-4163)
$ws2.Range("E2").Value = "2021-04-28T04:08:41.199009Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:18.614756"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/331/?format=json"

# Keep "data" as the active/visible tab (unchanged by the edit).
$ws1.Activate()
